$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "struggle"
$ws.Cells.Item(22, 3).Value = 4.677844420075353
$ws.Cells.Item(22, 4).Value = -3.651133604347696
$ws.Cells.Item(22, 5).Value = -7.842656075954431
$ws.Cells.Item(22, 6).Value = -0.7247915863990784
$ws.Cells.Item(22, 7).Value = -2.964529037475586
$ws.Cells.Item(22, 8).Value = -2.036930084228516
$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "struggle"
$ws.Cells.Item(23, 3).Value = -2.627654522657398
$ws.Cells.Item(23, 4).Value = -2.928949266672134
$ws.Cells.Item(23, 5).Value = 4.230176210403448
$ws.Cells.Item(23, 6).Value = -0.9755517840385436
$ws.Cells.Item(23, 7).Value = -3.013092756271362
$ws.Cells.Item(23, 8).Value = -0.1954768747091293
$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "struggle"
$ws.Cells.Item(24, 3).Value = -4.852406792342663
$ws.Cells.Item(24, 4).Value = 0.3913787733763447
$ws.Cells.Item(24, 5).Value = 0.2968738228082666
$ws.Cells.Item(24, 6).Value = -0.0148134818300604
$ws.Cells.Item(24, 7).Value = -4.330729007720947
$ws.Cells.Item(24, 8).Value = 0.6641632318496704
$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "struggle"
$ws.Cells.Item(25, 3).Value = -1.301035702228551
$ws.Cells.Item(25, 4).Value = 3.64691380783915
$ws.Cells.Item(25, 5).Value = -6.109266191720954
$ws.Cells.Item(25, 6).Value = -0.2535090744495392
$ws.Cells.Item(25, 7).Value = -4.50192403793335
$ws.Cells.Item(25, 8).Value = 0.8677340745925903
$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "struggle"
$ws.Cells.Item(26, 3).Value = 2.465943455696097
$ws.Cells.Item(26, 4).Value = -2.991184197366218
$ws.Cells.Item(26, 5).Value = -3.608212560415278
$ws.Cells.Item(26, 6).Value = -0.0308486949652433
$ws.Cells.Item(26, 7).Value = -3.680310487747192
$ws.Cells.Item(26, 8).Value = 1.009607553482056
$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "struggle"
$ws.Cells.Item(27, 3).Value = -1.307898223400096
$ws.Cells.Item(27, 4).Value = -2.068972408771528
$ws.Cells.Item(27, 5).Value = -0.7334359884262174
$ws.Cells.Item(27, 6).Value = -0.2924517393112182
$ws.Cells.Item(27, 7).Value = 0.6568328738212585
$ws.Cells.Item(27, 8).Value = 0.4216497242450714
$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "struggle"
$ws.Cells.Item(28, 3).Value = -1.702915767207749
$ws.Cells.Item(28, 4).Value = -0.5735956337302961
$ws.Cells.Item(28, 5).Value = -0.9715757742524092
$ws.Cells.Item(28, 6).Value = 0.683863639831543
$ws.Cells.Item(28, 7).Value = 4.383111000061035
$ws.Cells.Item(28, 8).Value = -1.505782842636108
$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "struggle"
$ws.Cells.Item(29, 3).Value = -3.414293382316824
$ws.Cells.Item(29, 4).Value = 0.2869436666369428
$ws.Cells.Item(29, 5).Value = -0.1008520126342796
$ws.Cells.Item(29, 6).Value = 0.6068946123123169
$ws.Cells.Item(29, 7).Value = 4.862334728240967
$ws.Cells.Item(29, 8).Value = -0.4990769028663635
$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "struggle"
$ws.Cells.Item(30, 3).Value = -3.79365611076355
$ws.Cells.Item(30, 4).Value = 1.782416181638838
$ws.Cells.Item(30, 5).Value = 1.588389292359353
$ws.Cells.Item(30, 6).Value = -0.4978551864624023
$ws.Cells.Item(30, 7).Value = 3.975052833557129
$ws.Cells.Item(30, 8).Value = 0.3481931984424591
$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "struggle"
$ws.Cells.Item(31, 3).Value = -3.179930865764618
$ws.Cells.Item(31, 4).Value = -0.5433011054992622
$ws.Cells.Item(31, 5).Value = 2.409818679094315
$ws.Cells.Item(31, 6).Value = -0.299934834241867
$ws.Cells.Item(31, 7).Value = 1.433395266532898
$ws.Cells.Item(31, 8).Value = 0.5711590051651001